# Insert a new weekly record (row) for "Camote" at row 12 of the
# "Hortaliza, Vega Modelo de Temuco - Camote" worksheet. This pushes the
# existing rows 12..126 down to rows 13..127, so the sheet dimension grows
# from A1:R126 to A1:R127, matching the weekly price-logic update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12, shifting rows below it down.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new market record.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44819
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 100114002
$ws.Range("G12").Value = "Camote"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 20000
$ws.Range("N12").Value = "$/malla 20 kilos"
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 1000
$ws.Range("Q12").Value = 20
$ws.Range("R12").Value = "Hortaliza"
